# The authored change swaps the contents of ppt/theme/theme1.xml ("Office
# Theme" colours) and ppt/theme/theme2.xml ("Integral" colours) so that the
# slide master's theme (theme2.xml, the design actually driving every
# slide) switches from the "Integral" palette to the default "Office"
# palette, while the notes-master theme (theme1.xml) ends up holding the
# "Integral" palette.  The font scheme and format scheme are identical
# between the two themes already, so only the 12 theme colours (and the
# theme/colour-scheme display names) actually change.
#
# The PowerPoint object model exposes the live theme colours through
# ThemeColorScheme (indices 1-12 = dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink), reachable from a Slide, the SlideMaster, the NotesMaster, etc.
# Here we drive it from the slide master so every slide picks up the new
# "Office" palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Target palette = the former theme1.xml ("Office Theme") colours.
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
